# Apply the acceptance-test-cases-g3 content update:
#  - Fill the "Actual Results" column (E) with the same PASS-evidence text
#    already present in the "Expected Result" column (D) for every
#    populated test-case row.
#  - Fill in the remaining blank "Pass/Fail" (F) cells with "PASS" for the
#    rows that didn't have it yet.
#  - Correct the AT301 test-data cell (C7) input format to match its
#    30/3/8K description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the shipment input format for AT301 (row 7) ---
$ws.Range("C7").Value = "SHIPMENT:`nWeight: 30kg`nSize: 3`nDestination: 8K`nInput Format: 30 3 8K"

# --- Rows whose Actual Results (E) should mirror Expected Result (D),
#     using the same centered/wrap style as column D on those rows ---
$centeredWrapRows = @(2, 4, 5, 7, 9, 13, 15, 17, 19)
foreach ($r in $centeredWrapRows) {
    $dst = $ws.Cells.Item($r, 5)
    $dst.Value = $ws.Cells.Item($r, 4).Text
    $dst.WrapText = $true
    $dst.VerticalAlignment = -4108
}

# --- Rows 11 & 12 use the plain wrap-only style (matching B/C/D there) ---
$wrapOnlyRows = @(11, 12)
foreach ($r in $wrapOnlyRows) {
    $dst = $ws.Cells.Item($r, 5)
    $dst.Value = $ws.Cells.Item($r, 4).Text
    $dst.WrapText = $true
}

# --- Fill in the missing "PASS" Pass/Fail cells ---
$passRows = @(11, 12, 13, 19)
foreach ($r in $passRows) {
    $dst = $ws.Cells.Item($r, 6)
    $dst.Value = "PASS"
    $dst.WrapText = $true
    $dst.VerticalAlignment = -4108
}
